$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New task text for week-3 rows (written in the order that reproduces the
#     original shared-string table ordering: 50..60) ---
$ws.Range("C59").Value = "activity diagrams maken"
$ws.Range("C60").Value = "Use case diagram afmaken"
$ws.Range("C61").Value = "Activity diagrams afmaken"
$ws.Range("C62").Value = "Beginnen aan het paper"
$ws.Range("C63").Value = "Code commenten"
$ws.Range("C65").Value = "Paper amaken"
$ws.Range("C66").Value = "Code een laatste keer checken"
$ws.Range("C67").Value = "Testen van de robot"
$ws.Range("C57").Value = "Code afmaken"
$ws.Range("C58").Value = "Code testen"
$ws.Range("C51").Value = "Werken aan de could haves"
$ws.Range("C54").Value = "Werken aan de could haves"

# --- Recolor the day blocks in column C for week 3 (rows 48-68), matching
#     the same yellow/red/green/blue banding used by the earlier weeks ---
$ws.Range("C48:C53").Interior.Color = 65535     # yellow - 3-apr
$ws.Range("C54:C59").Interior.Color = 255       # red    - 4-apr
$ws.Range("C60:C64").Interior.Color = 32768     # green  - 5-apr
$ws.Range("C65:C68").Interior.Color = 16737843  # blue   - 6-apr

# --- Update the saved view/selection state ---
$ws.Range("C54").Select() | Out-Null
